# Commit: fix(gui) step 1 and 2
# Updates the list date (A1) and the three price values (D33:D35)
# on the "Hoja1" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the printed date from 45308 (2024-01-17) to 45309 (2024-01-18)
$ws.Range("A1").Value = 45309

# Step 2: update the updated prices
$ws.Range("D33").Value = 1305
$ws.Range("D34").Value = 1275.478
$ws.Range("D35").Value = 949.728
